$d = $word.ActiveDocument

# 1. "Kicho E Jang" -> "Kicho I Jang"
$rng = $d.Content
$found = $rng.Find.Execute("Kicho E Jang", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $charRng = $d.Range($rng.Start + 6, $rng.Start + 7)
    $charRng.Text = "I"
}

# 2. "Kicho Sam Jan" -> "Kicho Sam Jang" (missing trailing "g")
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Kicho Sam Jan", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $lastCharRng = $d.Range($rng2.End - 1, $rng2.End)
    $lastCharRng.Text = "ng"
}

# 3. "Palgue E Jang" -> "Palgue I Jang"
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("Palgue E Jang", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $charRng3 = $d.Range($rng3.Start + 7, $rng3.Start + 8)
    $charRng3.Text = "I"
}
